$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.672.04'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '3.480.75'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '3.479.75'
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.141'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.00'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.430'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').Value = '4.075.63'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.41%  '
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '67.645.65'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000178'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('D18').Value = '3.478.67'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '395.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.540'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.94%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '72.09'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('E27').Value = '  +1.21%  '
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.177'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.39%  '
$ws.Range('E33').Value = '  +1.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.04'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.894'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.45%  '
$ws.Range('E41').Value = '  -1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.33'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0718'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.35'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.05%  '
$ws.Range('D47').Value = '2.748.29'
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '41.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0299'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '328.07'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.07%  '
